# Fruta / hortaliza, semanal
# Insert a new weekly record as row 182 in the "Membrillo" price sheet,
# pushing the former rows 182-193 down to 183-194.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 182 (shifts existing rows 182-193 -> 183-194)
$ws.Rows(182).Insert()

# Populate the newly inserted row 182 with the new weekly record
$ws.Range("A182").Value = 10
$ws.Range("B182").Value = "Vega Modelo de Temuco"
$ws.Range("C182").Value = "La Araucanía"
$ws.Range("D182").Value = 44753
$ws.Range("E182").Value = 9
$ws.Range("F182").Value = "Fruta"
$ws.Range("G182").Value = 100104
$ws.Range("H182").Value = "Frutos de pepita"
$ws.Range("I182").Value = 100104003
$ws.Range("J182").Value = "Membrillo"
$ws.Range("K182").Value = "Champion"
$ws.Range("L182").Value = "Primera"
$ws.Range("M182").Value = 25
$ws.Range("N182").Value = 10000
$ws.Range("O182").Value = 10000
$ws.Range("P182").Value = 10000
$ws.Range("Q182").Value = "$/bandeja 18 kilos granel"
$ws.Range("R182").Value = "Región de O'Higgins"
$ws.Range("S182").Value = 556
$ws.Range("T182").Value = 18
